$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (RUG357.fasta) - all remaining rows shift up by one
$ws.Rows.Item(2).Delete()

# Remove the "max" column (old column C) - prediction/rejection-f shift left
$ws.Columns.Item(3).Delete()

# Update the "1-f__UBA660" numeric values and the rejection-f text for each row
$ws.Range("B2").Value = 3450.486807122075
$ws.Range("B3").Value = 11029.99761819572
$ws.Range("B4").Value = 3551.717772586602
$ws.Range("B5").Value = -1271.722072691366
$ws.Range("B6").Value = 6385.3487251236
$ws.Range("B7").Value = 18130.80467507385
$ws.Range("B8").Value = 21435.34895983516
$ws.Range("B9").Value = 5725.519969701298
$ws.Range("B10").Value = 19170.29751393408
$ws.Range("B11").Value = 28457.46456298053

# RUG705.fasta's prediction is now rejected because its value is negative
$ws.Range("D5").Value = "f__UBA660(reject)"
